# "Adatfájl változtatás+programok fixálása az első teszt után"
#
# 1. "kezdőhelyek" sheet: selection moves from H13 to C10
# 2. "jatekos szinek" (player colours) sheet becomes the active tab;
#    player 8's colour changes from "lavender" to "maroon" (this also
#    prunes the now-unused "lavender" shared string and appends a new
#    "maroon" shared string, automatically reindexing every other sheet
#    that references the resource-name strings after it)
# 3. "terkep" (map) sheet: loses the active-tab flag; the stray
#    harvester-count value of 12 in D5:D7 is cleared (D4/D8 stay as-is)
# 4. Each of the 12 per-player resource sheets ("1".."12") had its
#    quantity column (B2:B7) filled in for testing; that test data is
#    cleared back out again. Their selections are also updated.

$wb = $excel.ActiveWorkbook

# --- sheet "kezdőhelyek" ---------------------------------------------------
$wsStart = $wb.Worksheets.Item(1)
$wsStart.Activate()
$wsStart.Range("C10").Select()

# --- per-player resource sheets "1".."12" (workbook order: items 4-15) ----
# Clear the test quantities that were typed into column B.
$wb.Worksheets.Item(4).Range("B2:B7").ClearContents()   # sheet "1"
$wb.Worksheets.Item(5).Range("B2:B7").ClearContents()   # sheet "2"
$wb.Worksheets.Item(6).Range("B2:B7").ClearContents()   # sheet "3"
$wb.Worksheets.Item(7).Range("B2:B7").ClearContents()   # sheet "4"
$wb.Worksheets.Item(8).Range("B2:B7").ClearContents()   # sheet "5"
$wb.Worksheets.Item(10).Range("B2:B7").ClearContents()  # sheet "7"

# Update the remembered selections on each of these sheets.
$wb.Worksheets.Item(4).Range("E14").Select()    # sheet "1" -> E14
$wb.Worksheets.Item(8).Range("B2:B7").Select()  # sheet "5" -> B2:B7
$wb.Worksheets.Item(10).Range("B2").Select()    # sheet "7" -> B2

# --- sheet "térkép" ---------------------------------------------------------
$wsMap = $wb.Worksheets.Item(3)
$wsMap.Range("D5:D7").ClearContents()
$wsMap.Range("D6").Select()

# --- sheet "játékos színek" (player colours) -------------------------------
$wsColors = $wb.Worksheets.Item(2)
$wsColors.Range("B10").Value = "maroon"

# Make this the active sheet/tab (was "térkép" before).
$wsColors.Activate()
$wsColors.Range("G7").Select()
